# Fruta / hortaliza, semanal
# The weekly data rows for "Comercializadora del Agro de Limarí - Arándano (blue)"
# got reshuffled: the Fecha/Volumen/Precio columns (D, M, N, O, P, S) for each
# data row (2-18, row 5 unaffected) now hold the values that used to belong to a
# different row. Capture all current values first, then write the permuted
# values back so that reads never see already-overwritten data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "M", "N", "O", "P", "S")

# destination row -> source row (values currently sitting in source row move to destination row)
$mapping = @{
    2  = 10
    3  = 17
    4  = 18
    6  = 9
    7  = 3
    8  = 14
    9  = 12
    10 = 11
    11 = 4
    12 = 6
    13 = 16
    14 = 8
    15 = 13
    16 = 15
    17 = 2
    18 = 7
}

# Snapshot current values for every row referenced by the mapping.
$snapshot = @{}
foreach ($row in $mapping.Values) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range($col + $row).Value2
    }
    $snapshot[$row] = $rowVals
}

# Write the permuted values into their destination rows.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $rowVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range($col + $destRow).Value = $rowVals[$col]
    }
}
